$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("c2-c4")

# Populate the new "GWP-LULUC_eol" column (L) with 0 for all data rows that
# are missing it (rows 2-284). A handful of these cells already exist as
# empty-but-styled cells (date-format style carried over from column K/F);
# normalize their style back to Normal to match the recalculated output,
# then fill in the value.
$ws.Range("L2:L284").Style = "Normal"
$ws.Range("L2:L284").Value = 0

# The last few rows (285-288) were missing several trailing columns
# entirely (G, K, L) - fill them in with 0 as well.
$ws.Range("G285").Value = 0
$ws.Range("K285:K288").Value = 0
$ws.Range("L285:L288").Value = 0

# Update the sheet's active selection/view to match where the edits were
# made.
$ws.Select()
$ws.Range("G285").Select()
